$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 39140
$ws.Range("J87").Value = 39140
$ws.Range("L87").Value = 39140
$ws.Range("N87").Value = -41636

$ws.Range("H90").Value = 39140
$ws.Range("J90").Value = 39140
$ws.Range("L90").Value = 117420
$ws.Range("N90").Value = -129900

$ws.Range("H93").Value = 50101
$ws.Range("J93").Value = 50101
$ws.Range("L93").Value = 50101
$ws.Range("N93").Value = -55093

$ws.Range("H113").Value = 2975.4614
$ws.Range("I113").Value = 2148.3333
$ws.Range("J113").Value = 4103.364
$ws.Range("K113").Value = 2148.3333
$ws.Range("L113").Value = 4103.364
$ws.Range("M113").Value = 1105.6667
$ws.Range("N113").Value = -10611.364

$ws.Range("H138").Value = 1863.129
$ws.Range("I138").Value = 1084.619
$ws.Range("J138").Value = 3498
$ws.Range("K138").Value = 3253.857
$ws.Range("L138").Value = 10494
$ws.Range("M138").Value = 1886.143
$ws.Range("N138").Value = -20774

$ws.Range("H141").Value = 4936.3774
$ws.Range("I141").Value = 1410.6666
$ws.Range("J141").Value = 12402.588
$ws.Range("K141").Value = 4231.9998
$ws.Range("L141").Value = 37207.764
$ws.Range("M141").Value = 948.0002000000004
$ws.Range("N141").Value = -47567.764

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3106.42
$ws.Range("I32").Value = 2733.932
$ws.Range("J32").Value = 5838
$ws.Range("K32").Value = 2733.932
$ws.Range("L32").Value = 5838
$ws.Range("M32").Value = -2446.932
$ws.Range("N32").Value = -6412

$ws.Range("H61").Value = 1758.8379
$ws.Range("I61").Value = 1571.0869
$ws.Range("J61").Value = 2067.2856
$ws.Range("K61").Value = 1571.0869
$ws.Range("L61").Value = 2067.2856
$ws.Range("M61").Value = -1359.0869
$ws.Range("N61").Value = -2491.2856

$ws.Range("H63").Value = 3520.606
$ws.Range("J63").Value = 9093.333000000001
$ws.Range("L63").Value = 9093.333000000001
$ws.Range("N63").Value = -10465.333

$ws.Range("H66").Value = 3520.606
$ws.Range("J66").Value = 9093.333000000001
$ws.Range("L66").Value = 45466.665
$ws.Range("N66").Value = -52330.665

$ws.Range("H74").Value = 1289.942
$ws.Range("I74").Value = 915.7292
$ws.Range("J74").Value = 2145.2856
$ws.Range("K74").Value = 915.7292
$ws.Range("L74").Value = 2145.2856
$ws.Range("M74").Value = -41.72919999999999
$ws.Range("N74").Value = -3893.2856

$ws.Range("H77").Value = 1289.942
$ws.Range("I77").Value = 915.7292
$ws.Range("J77").Value = 2145.2856
$ws.Range("K77").Value = 4578.646
$ws.Range("L77").Value = 10726.428
$ws.Range("M77").Value = -210.6459999999997
$ws.Range("N77").Value = -19462.428

$ws.Range("H132").Value = 2019949.8
$ws.Range("I132").Value = 2193.0286
$ws.Range("J132").Value = 4635560
$ws.Range("K132").Value = 6579.085800000001
$ws.Range("L132").Value = 13906680
$ws.Range("M132").Value = -4049.085800000001
$ws.Range("N132").Value = -13911740

$ws.Range("H136").Value = 1758.8379
$ws.Range("I136").Value = 1571.0869
$ws.Range("J136").Value = 2067.2856
$ws.Range("K136").Value = 4713.2607
$ws.Range("L136").Value = 6201.8568
$ws.Range("M136").Value = -2163.2607
$ws.Range("N136").Value = -11301.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2241.7273
$ws.Range("I105").Value = 2079.8
$ws.Range("J105").Value = 2376.6667
$ws.Range("K105").Value = 2079.8
$ws.Range("L105").Value = 2376.6667
$ws.Range("M105").Value = -332.8000000000002
$ws.Range("N105").Value = -5870.6667

$ws.Range("H132").Value = 39315
$ws.Range("J132").Value = 39315
$ws.Range("L132").Value = 39315
$ws.Range("N132").Value = -49435

$ws.Range("H134").Value = 3606.611
$ws.Range("I134").Value = 1393.9762
$ws.Range("J134").Value = 6704.3
$ws.Range("K134").Value = 4181.9286
$ws.Range("L134").Value = 20112.9
$ws.Range("M134").Value = -1646.9286
$ws.Range("N134").Value = -25182.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8132459.5
$ws.Range("I31").Value = 2140.68
$ws.Range("J31").Value = 20836084
$ws.Range("K31").Value = 2140.68
$ws.Range("L31").Value = 20836084
$ws.Range("M31").Value = -1845.68
$ws.Range("N31").Value = -20836674

$ws.Range("H34").Value = 8132459.5
$ws.Range("I34").Value = 2140.68
$ws.Range("J34").Value = 20836084
$ws.Range("K34").Value = 2140.68
$ws.Range("L34").Value = 20836084
$ws.Range("M34").Value = -1938.68
$ws.Range("N34").Value = -20836488

$ws.Range("H99").Value = 3264.0476
$ws.Range("I99").Value = 3126.5386
$ws.Range("J99").Value = 3487.5
$ws.Range("K99").Value = 3126.5386
$ws.Range("L99").Value = 3487.5
$ws.Range("M99").Value = -1628.5386
$ws.Range("N99").Value = -6483.5

$ws.Range("H126").Value = 3264.0476
$ws.Range("I126").Value = 3126.5386
$ws.Range("J126").Value = 3487.5
$ws.Range("K126").Value = 9379.6158
$ws.Range("L126").Value = 10462.5
$ws.Range("M126").Value = -6909.6158
$ws.Range("N126").Value = -15402.5

$ws.Range("H132").Value = 2089.7222
$ws.Range("I132").Value = 1293.9667
$ws.Range("J132").Value = 3084.4167
$ws.Range("K132").Value = 3881.9001
$ws.Range("L132").Value = 9253.250100000001
$ws.Range("M132").Value = -1351.9001
$ws.Range("N132").Value = -14313.2501

$ws.Range("H134").Value = 1722.6285
$ws.Range("I134").Value = 933.25
$ws.Range("J134").Value = 3444.9092
$ws.Range("K134").Value = 2799.75
$ws.Range("L134").Value = 10334.7276
$ws.Range("M134").Value = -264.75
$ws.Range("N134").Value = -15404.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 85
$ws.Range("I2").Value = 57.333332
$ws.Range("J2").Value = 101.6
$ws.Range("K2").Value = 343.999992
$ws.Range("L2").Value = 609.5999999999999
$ws.Range("M2").Value = -230.999992
$ws.Range("N2").Value = -835.5999999999999

$ws.Range("H12").Value = 353.375
$ws.Range("I12").Value = 200
$ws.Range("J12").Value = 375.2857
$ws.Range("K12").Value = 600
$ws.Range("L12").Value = 1125.8571
$ws.Range("M12").Value = -427
$ws.Range("N12").Value = -1471.8571

$ws.Range("H80").Value = 1343.75
$ws.Range("J80").Value = 1606.1428
$ws.Range("L80").Value = 4818.428400000001
$ws.Range("N80").Value = -6690.428400000001

$ws.Range("H83").Value = 1343.75
$ws.Range("J83").Value = 1606.1428
$ws.Range("L83").Value = 14455.2852
$ws.Range("N83").Value = -23815.2852

$ws.Range("H122").Value = 2657.3774
$ws.Range("I122").Value = 1070.6
$ws.Range("J122").Value = 3026.3953
$ws.Range("K122").Value = 9635.4
$ws.Range("L122").Value = 27237.5577
$ws.Range("M122").Value = -7185.4
$ws.Range("N122").Value = -32137.5577

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 117187.5
$ws.Range("J135").Value = 117187.5
$ws.Range("L135").Value = 117187.5
$ws.Range("N135").Value = -127327.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1832.1538
$ws.Range("I16").Value = 897.0909
$ws.Range("J16").Value = 6975
$ws.Range("K16").Value = 897.0909
$ws.Range("L16").Value = 6975
$ws.Range("M16").Value = -727.0909
$ws.Range("N16").Value = -7315

$ws.Range("H40").Value = 33336324
$ws.Range("I40").Value = 40002252
$ws.Range("J40").Value = 6680
$ws.Range("K40").Value = 40002252
$ws.Range("L40").Value = 6680
$ws.Range("M40").Value = -40002116
$ws.Range("N40").Value = -6952

$ws.Range("H55").Value = 428.6154
$ws.Range("I55").Value = 388
$ws.Range("J55").Value = 564
$ws.Range("K55").Value = 388
$ws.Range("L55").Value = 564
$ws.Range("M55").Value = -215
$ws.Range("N55").Value = -910

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 16449
$ws.Range("J92").Value = 16449
$ws.Range("L92").Value = 16449
$ws.Range("N92").Value = -21441

$ws.Range("H122").Value = 4846.4414
$ws.Range("I122").Value = 3448.7083
$ws.Range("K122").Value = 10346.1249
$ws.Range("M122").Value = -7896.124899999999
